{"js": "// Collapse each \"<id>p083r_aN</id>\" (currently split across three runs:\n// the literal \"<id>\" tag, the bare id text, and the literal \"</id>\" tag)\n// into a single run containing \"<id>p083r_N</id>\" \u2014 i.e. strip the \"a\"\n// out of the id value while merging the three runs into one.\nconst ids = [\n  [\"p083r_a1\", \"p083r_1\"],\n  [\"p083r_a2\", \"p083r_2\"],\n  [\"p083r_a3\", \"p083r_3\"],\n  [\"p083r_a4\", \"p083r_4\"],\n  [\"p083r_a5\", \"p083r_5\"],\n  [\"p083r_a6\", \"p083r_6\"],\n];\n\nfor (const [oldId, newId] of ids) {\n  const oldText = \"<id>\" + oldId + \"</id>\";\n  const newText = \"<id>\" + newId + \"</id>\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Collapse each \"<id>p083r_aN</id>\" (currently split across three runs:\n# the literal \"<id>\" tag, the bare id text, and the literal \"</id>\" tag)\n# into a single run containing \"<id>p083r_N</id>\" - i.e. strip the \"a\"\n# out of the id value while merging the three runs into one.\n$d = $word.ActiveDocument\n\n$ids = @(\n  @(\"p083r_a1\", \"p083r_1\"),\n  @(\"p083r_a2\", \"p083r_2\"),\n  @(\"p083r_a3\", \"p083r_3\"),\n  @(\"p083r_a4\", \"p083r_4\"),\n  @(\"p083r_a5\", \"p083r_5\"),\n  @(\"p083r_a6\", \"p083r_6\")\n)\n\nforeach ($pair in $ids) {\n    $oldText = \"<id>\" + $pair[0] + \"</id>\"\n    $newText = \"<id>\" + $pair[1] + \"</id>\"\n\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
